# Scheduled market-price refresh for the Leve profit-tracking sheets.
# Updates currentAveragePrice / NQ / HQ / LevePrice(NQ|HQ) / LeveProfit(NQ|HQ)
# columns (H:N) with the latest scraped values for the affected leves.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3: Leather Grimoire
$ws.Range("H3").Value = 30657
$ws.Range("J3").Value = 30657
$ws.Range("L3").Value = 30657
$ws.Range("N3").Value = -30885
# Row 64: Void Glue
$ws.Range("H64").Value = 5286.4375
$ws.Range("I64").Value = 4041.111
$ws.Range("J64").Value = 6887.5713
$ws.Range("K64").Value = 4041.111
$ws.Range("L64").Value = 6887.5713
$ws.Range("M64").Value = -3793.111
$ws.Range("N64").Value = -7383.5713
# Row 67: Void Glue
$ws.Range("H67").Value = 5286.4375
$ws.Range("I67").Value = 4041.111
$ws.Range("J67").Value = 6887.5713
$ws.Range("K67").Value = 4041.111
$ws.Range("L67").Value = 6887.5713
$ws.Range("M67").Value = -3183.111
$ws.Range("N67").Value = -8603.5713
# Row 86: Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 1810.3529
$ws.Range("I86").Value = 1188.2222
$ws.Range("J86").Value = 2510.25
$ws.Range("K86").Value = 1188.2222
$ws.Range("L86").Value = 2510.25
$ws.Range("M86").Value = -65.22219999999993
$ws.Range("N86").Value = -4756.25
# Row 87: Noble Gold
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
# Row 89: Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 1810.3529
$ws.Range("I89").Value = 1188.2222
$ws.Range("J89").Value = 2510.25
$ws.Range("K89").Value = 5941.111
$ws.Range("L89").Value = 12551.25
$ws.Range("M89").Value = -325.1109999999999
$ws.Range("N89").Value = -23783.25
# Row 90: Noble Gold
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
# Row 92: Enchanted Koppranickel Ink
$ws.Range("H92").Value = 927204.0600000001
$ws.Range("I92").Value = 1112594.9
$ws.Range("K92").Value = 1112594.9
$ws.Range("M92").Value = -1111346.9
# Row 102: Marid Leather Grimoire
$ws.Range("H102").Value = 30657
$ws.Range("J102").Value = 30657
$ws.Range("L102").Value = 30657
$ws.Range("N102").Value = -37147
# Row 116: Growth Formula Kappa
$ws.Range("H116").Value = 2844.5454
$ws.Range("J116").Value = 1500
$ws.Range("L116").Value = 1500
$ws.Range("N116").Value = -8384
# Row 132: Growth Formula Lambda
$ws.Range("H132").Value = 46049.78
$ws.Range("I132").Value = 46049.78
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 138149.34
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -135619.34
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 18569.29
$ws.Range("I32").Value = 4004.6094
$ws.Range("J32").Value = 67629.266
$ws.Range("K32").Value = 4004.6094
$ws.Range("L32").Value = 67629.266
$ws.Range("M32").Value = -3717.6094
$ws.Range("N32").Value = -68203.266
# Row 63: Mythrite Rivets
$ws.Range("H63").Value = 8236.909
$ws.Range("I63").Value = 9974.375
$ws.Range("J63").Value = 3603.6667
$ws.Range("K63").Value = 9974.375
$ws.Range("L63").Value = 3603.6667
$ws.Range("M63").Value = -9288.375
$ws.Range("N63").Value = -4975.6667
# Row 66: Mythrite Rivets
$ws.Range("H66").Value = 8236.909
$ws.Range("I66").Value = 9974.375
$ws.Range("J66").Value = 3603.6667
$ws.Range("K66").Value = 49871.875
$ws.Range("L66").Value = 18018.3335
$ws.Range("M66").Value = -46439.875
$ws.Range("N66").Value = -24882.3335
# Row 102: Tama-hagane Ingot
$ws.Range("H102").Value = 4870.3335
$ws.Range("I102").Value = 4870.3335
$ws.Range("K102").Value = 4870.3335
$ws.Range("M102").Value = -3248.3335
# Row 122: High Durium Nugget
$ws.Range("H122").Value = 2661.1765
$ws.Range("I122").Value = 2074.1667
$ws.Range("J122").Value = 4070
$ws.Range("K122").Value = 6222.500100000001
$ws.Range("L122").Value = 12210
$ws.Range("M122").Value = -3772.500100000001
$ws.Range("N122").Value = -17110

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Oroshigane Ingot
$ws.Range("H99").Value = 2323
$ws.Range("I99").Value = 2149.875
$ws.Range("J99").Value = 2600
$ws.Range("K99").Value = 2149.875
$ws.Range("L99").Value = 2600
$ws.Range("M99").Value = -651.875
$ws.Range("N99").Value = -5596
# Row 134: Ruthenium Ingot
$ws.Range("H134").Value = 3823.3684
$ws.Range("I134").Value = 2136.2693
$ws.Range("K134").Value = 6408.8079
$ws.Range("M134").Value = -3873.8079

$ws = $wb.Worksheets.Item("CRP")
# Row 36: Steel Spear
$ws.Range("H36").Value = 2166.3333
$ws.Range("I36").Value = 2166.3333
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2166.3333
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1778.3333
$ws.Range("N36").ClearContents()
# Row 40: Steel Spear
$ws.Range("H40").Value = 2166.3333
$ws.Range("I40").Value = 2166.3333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2166.3333
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2006.3333
$ws.Range("N40").ClearContents()
# Row 107: White Oak Lumber
$ws.Range("H107").Value = 414.45456
$ws.Range("I107").Value = 339.77777
$ws.Range("K107").Value = 339.77777
$ws.Range("M107").Value = 1580.22223
# Row 134: Ceiba Lumber
$ws.Range("H134").Value = 3869.0952
$ws.Range("I134").Value = 1716.909
$ws.Range("J134").Value = 6236.5
$ws.Range("K134").Value = 5150.727000000001
$ws.Range("L134").Value = 18709.5
$ws.Range("M134").Value = -2615.727000000001
$ws.Range("N134").Value = -23779.5

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Grilled Trout
$ws.Range("H3").Value = 6543.3335
$ws.Range("I3").Value = 1696.6666
$ws.Range("K3").Value = 5089.9998
$ws.Range("M3").Value = -4977.9998
# Row 21: Raw Oyster
$ws.Range("H21").Value = 210.4
$ws.Range("I21").Value = 200.5
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 601.5
$ws.Range("L21").Value = 750
$ws.Range("M21").Value = -428.5
$ws.Range("N21").Value = -1096
# Row 113: Night Vinegar
$ws.Range("H113").Value = 504.91666
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 511.08823
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 1533.26469
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -5873.26469
# Row 131: Tsai tou Vounou
$ws.Range("H131").Value = 8773644
$ws.Range("J131").Value = 9261038
$ws.Range("L131").Value = 27783114
$ws.Range("N131").Value = -27793194

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Hardsilver Ingot
$ws.Range("H80").Value = 2971.7
$ws.Range("I80").Value = 2892.8572
$ws.Range("J80").Value = 3155.6667
$ws.Range("K80").Value = 2892.8572
$ws.Range("L80").Value = 3155.6667
$ws.Range("M80").Value = -1894.8572
$ws.Range("N80").Value = -5151.6667
# Row 83: Hardsilver Ingot
$ws.Range("H83").Value = 2971.7
$ws.Range("I83").Value = 2892.8572
$ws.Range("J83").Value = 3155.6667
$ws.Range("K83").Value = 14464.286
$ws.Range("L83").Value = 15778.3335
$ws.Range("M83").Value = -9472.286
$ws.Range("N83").Value = -25762.3335
# Row 133: Lar Pendulums
$ws.Range("H133").Value = 47262.5
$ws.Range("J133").Value = 47262.5
$ws.Range("L133").Value = 47262.5
$ws.Range("N133").Value = -57382.5
# Row 138: White Gold Halfmask of Maiming
$ws.Range("H138").Value = 67378.5
$ws.Range("J138").Value = 67378.5
$ws.Range("L138").Value = 67378.5
$ws.Range("N138").Value = -77658.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Aldgoat Leather
$ws.Range("H22").Value = 7528.125
$ws.Range("I22").Value = 1050
$ws.Range("J22").Value = 10472.728
$ws.Range("K22").Value = 1050
$ws.Range("L22").Value = 10472.728
$ws.Range("M22").Value = -755
$ws.Range("N22").Value = -11062.728
# Row 27: Aldgoat Leather
$ws.Range("H27").Value = 7528.125
$ws.Range("I27").Value = 1050
$ws.Range("J27").Value = 10472.728
$ws.Range("K27").Value = 1050
$ws.Range("L27").Value = 10472.728
$ws.Range("M27").Value = -943
$ws.Range("N27").Value = -10686.728
# Row 45: Peisteskin Crakows
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# Row 46: Boar Leather
$ws.Range("H46").Value = 2116.6667
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 2425
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 2425
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -2801
# Row 68: Wyvern Leather
$ws.Range("H68").Value = 1869
$ws.Range("I68").Value = 1861.25
$ws.Range("J68").Value = 1900
$ws.Range("K68").Value = 1861.25
$ws.Range("L68").Value = 1900
$ws.Range("M68").Value = -1112.25
$ws.Range("N68").Value = -3398
# Row 71: Wyvern Leather
$ws.Range("H71").Value = 1869
$ws.Range("I71").Value = 1861.25
$ws.Range("J71").Value = 1900
$ws.Range("K71").Value = 9306.25
$ws.Range("L71").Value = 9500
$ws.Range("M71").Value = -5562.25
$ws.Range("N71").Value = -16988
# Row 75: Dhalmelskin Leggings of Aiming
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78: Dhalmelskin Leggings of Aiming
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 100: Tiger Leather
$ws.Range("H100").Value = 2086171.9
$ws.Range("I100").Value = 5684060
$ws.Range("K100").Value = 5684060
$ws.Range("M100").Value = -5683519
# Row 122: Gaja Leather
$ws.Range("H122").Value = 3214.4075
$ws.Range("I122").Value = 2012.375
$ws.Range("J122").Value = 3720.5264
$ws.Range("K122").Value = 6037.125
$ws.Range("L122").Value = 11161.5792
$ws.Range("M122").Value = -3587.125
$ws.Range("N122").Value = -16061.5792

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Rainbow Cloth
$ws.Range("H62").Value = 19105.666
$ws.Range("I62").Value = 22000.223
$ws.Range("K62").Value = 22000.223
$ws.Range("M62").Value = -21376.223
# Row 65: Rainbow Cloth
$ws.Range("H65").Value = 19105.666
$ws.Range("I65").Value = 22000.223
$ws.Range("K65").Value = 110001.115
$ws.Range("M65").Value = -106881.115
# Row 100: Kudzu Thread
$ws.Range("H100").Value = 700.5
$ws.Range("I100").Value = 699.5
$ws.Range("K100").Value = 1399
$ws.Range("M100").Value = -858
# Row 126: Snow Linen
$ws.Range("H126").Value = 72653.92999999999
$ws.Range("J126").Value = 1662.5
$ws.Range("L126").Value = 4987.5
$ws.Range("N126").Value = -9927.5
# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 2210.0789
$ws.Range("I132").Value = 2033.8889
$ws.Range("J132").Value = 2642.5454
$ws.Range("K132").Value = 6101.6667
$ws.Range("L132").Value = 7927.6362
$ws.Range("M132").Value = -3571.6667
$ws.Range("N132").Value = -12987.6362
